$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 20).Value = 1.74
$ws.Cells.Item(3, 7).Value = 2.5
$ws.Cells.Item(3, 9).Value = 5.2
$ws.Cells.Item(3, 12).Value = 1.5
$ws.Cells.Item(3, 17).Value = 2.68
$ws.Cells.Item(3, 20).Value = 2.24
$ws.Cells.Item(3, 21).Value = 1.64
$ws.Cells.Item(3, 25).Value = 12.5
$ws.Cells.Item(4, 6).Value = 1.09
$ws.Cells.Item(4, 7).Value = 2.62
$ws.Cells.Item(4, 8).Value = 2.74
$ws.Cells.Item(4, 9).Value = 980
$ws.Cells.Item(4, 10).Value = 3.6
$ws.Cells.Item(4, 11).Value = 980
$ws.Cells.Item(4, 13).Value = 1.03
$ws.Cells.Item(4, 14).Value = 1.1
$ws.Cells.Item(4, 16).Value = 2.18
$ws.Cells.Item(4, 17).Value = 1.45
$ws.Cells.Item(4, 18).Value = 1.53
$ws.Cells.Item(4, 22).Value = 1.39
$ws.Cells.Item(4, 23).Value = 1.62
$ws.Cells.Item(5, 12).Value = 1.65
$ws.Cells.Item(5, 13).Value = 1.15
$ws.Cells.Item(5, 14).Value = 2.26
$ws.Cells.Item(5, 15).Value = 1.64
$ws.Cells.Item(5, 16).Value = 1.41
$ws.Cells.Item(5, 18).Value = 1.14
$ws.Cells.Item(5, 19).Value = 6.4
$ws.Cells.Item(5, 20).Value = 2.24
$ws.Cells.Item(5, 21).Value = 1.66
$ws.Cells.Item(5, 24).Value = 7.4
$ws.Cells.Item(5, 25).Value = 8.800000000000001
$ws.Cells.Item(5, 26).Value = 22
$ws.Cells.Item(5, 27).Value = 80
$ws.Cells.Item(5, 28).Value = 7.6
$ws.Cells.Item(5, 29).Value = 7.2
$ws.Cells.Item(5, 30).Value = 17
$ws.Cells.Item(5, 31).Value = 65
$ws.Cells.Item(5, 32).Value = 16.5
$ws.Cells.Item(5, 33).Value = 14.5
$ws.Cells.Item(5, 34).Value = 28
$ws.Cells.Item(5, 35).Value = 120
$ws.Cells.Item(5, 36).Value = 55
$ws.Cells.Item(5, 37).Value = 980
$ws.Cells.Item(5, 38).Value = 110
$ws.Cells.Item(5, 40).Value = 75
$ws.Cells.Item(5, 41).Value = 120
$ws.Cells.Item(6, 6).Value = 1.16
$ws.Cells.Item(6, 7).Value = 1.71
$ws.Cells.Item(6, 8).Value = 1.09
$ws.Cells.Item(6, 9).Value = 9.6
$ws.Cells.Item(6, 10).Value = 4.2
$ws.Cells.Item(6, 11).Value = 980
$ws.Cells.Item(6, 12).Value = 1.01
$ws.Cells.Item(6, 13).Value = 1.01
$ws.Cells.Item(6, 14).Value = 2.12
$ws.Cells.Item(6, 15).Value = 1.25
$ws.Cells.Item(6, 16).Value = 2.12
$ws.Cells.Item(6, 18).Value = 1.37
$ws.Cells.Item(6, 19).Value = 2.44
$ws.Cells.Item(6, 20).Value = 1.01
$ws.Cells.Item(6, 21).Value = 1.01
$ws.Cells.Item(6, 22).Value = 1.11
$ws.Cells.Item(6, 23).Value = 2.4
$ws.Cells.Item(6, 24).Value = 1000
$ws.Cells.Item(6, 25).Value = 1000
$ws.Cells.Item(6, 26).Value = 1000
$ws.Cells.Item(6, 27).Value = 1000
$ws.Cells.Item(6, 28).Value = 1000
$ws.Cells.Item(6, 29).Value = 1000
$ws.Cells.Item(6, 30).Value = 1000
$ws.Cells.Item(6, 31).Value = 1000
$ws.Cells.Item(6, 32).Value = 1000
$ws.Cells.Item(6, 33).Value = 1000
$ws.Cells.Item(6, 34).Value = 1000
$ws.Cells.Item(6, 35).Value = 1000
$ws.Cells.Item(6, 36).Value = 1000
$ws.Cells.Item(6, 37).Value = 1000
$ws.Cells.Item(6, 38).Value = 1000
$ws.Cells.Item(6, 39).Value = 1000
$ws.Cells.Item(6, 40).Value = 1000
$ws.Cells.Item(6, 41).Value = 1000
$ws.Cells.Item(7, 6).Value = 2.46
$ws.Cells.Item(7, 12).Value = 1.01
$ws.Cells.Item(7, 13).Value = 1.06
$ws.Cells.Item(7, 14).Value = 1.01
$ws.Cells.Item(7, 15).Value = 1.36
$ws.Cells.Item(7, 17).Value = 2
$ws.Cells.Item(7, 18).Value = 1.1
$ws.Cells.Item(7, 19).Value = 1.01
$ws.Cells.Item(7, 20).Value = 1.01
$ws.Cells.Item(7, 21).Value = 1.01
$ws.Cells.Item(7, 22).Value = 1.39
$ws.Cells.Item(7, 23).Value = 1.55
$ws.Cells.Item(7, 24).Value = 19
$ws.Cells.Item(7, 25).Value = 17
$ws.Cells.Item(7, 26).Value = 30
$ws.Cells.Item(7, 27).Value = 75
$ws.Cells.Item(7, 28).Value = 15
$ws.Cells.Item(7, 29).Value = 11
$ws.Cells.Item(7, 30).Value = 19
$ws.Cells.Item(7, 31).Value = 50
$ws.Cells.Item(7, 32).Value = 25
$ws.Cells.Item(7, 33).Value = 18
$ws.Cells.Item(7, 34).Value = 26
$ws.Cells.Item(7, 35).Value = 70
$ws.Cells.Item(7, 36).Value = 55
$ws.Cells.Item(7, 37).Value = 44
$ws.Cells.Item(7, 38).Value = 65
$ws.Cells.Item(7, 39).Value = 1000
$ws.Cells.Item(7, 40).Value = 1000
$ws.Cells.Item(7, 41).Value = 1000
$ws.Cells.Item(9, 24).Value = 10.5
$ws.Cells.Item(9, 26).Value = 70
$ws.Cells.Item(9, 27).Value = 340
$ws.Cells.Item(9, 35).Value = 180
$ws.Cells.Item(9, 36).Value = 19
$ws.Cells.Item(10, 6).Value = 1.13
$ws.Cells.Item(10, 9).Value = 1000
$ws.Cells.Item(10, 16).Value = 2.96
